# Append 4 new order rows (rows 6-9) to the order history sheet,
# mirroring the structure/columns of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        A = "20250308000711"
        B = "2025-03-08 00:07:11"
        C = "[{'Item Name': 'Veg Thali', 'Price': 150.0}]"
        D = 150
        E = "Paid"
        F = "Take Away"
        G = "Pending"
        H = "nan"
    },
    @{
        A = "20250308000718"
        B = "2025-03-08 00:07:18"
        C = "[{'Item Name': 'Dal Rice', 'Price': 120.0}]"
        D = 120
        E = "Paid"
        F = "Take Away"
        G = "Delivered"
        H = "nan"
    },
    @{
        A = "20250308155356"
        B = "2025-03-08 15:53:56"
        C = "[{'Item Name': 'South Indian Thali', 'Price': 200.0}, {'Item Name': 'South Indian Thali', 'Price': 200.0}]"
        D = 400
        E = "Not Paid"
        F = "Take Away"
        G = "Pending"
        H = "nan"
    },
    @{
        A = "20250308155419"
        B = "2025-03-08 15:54:19"
        C = "[{'Item Name': 'Veg Thali', 'Price': 150.0}, {'Item Name': 'Roti Sabzi', 'Price': 100.0}, {'Item Name': 'Roti Sabzi', 'Price': 100.0}, {'Item Name': 'Roti Sabzi', 'Price': 100.0}]"
        D = 450
        E = "Not Paid"
        F = "Take Away"
        G = "Pending"
        H = "rajas 1"
    }
)

$startRow = 6
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds a long numeric-looking order id that must stay text,
    # like the rest of the sheet (inline/shared string, not a number).
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row.A
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
}
